$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 461
$ws.Cells.Item(12, 9).Value = 201
$ws.Cells.Item(12, 10).Value = 526
$ws.Cells.Item(12, 11).Value = 201
$ws.Cells.Item(12, 12).Value = 526
$ws.Cells.Item(12, 13).Value = -31
$ws.Cells.Item(12, 14).Value = -866
$ws.Cells.Item(19, 8).Value = 1386.1
$ws.Cells.Item(19, 10).Value = 1193
$ws.Cells.Item(19, 12).Value = 1193
$ws.Cells.Item(19, 14).Value = -1543
$ws.Cells.Item(55, 8).Value = 275
$ws.Cells.Item(55, 9).Value = 275
$ws.Cells.Item(55, 11).Value = 275
$ws.Cells.Item(55, 13).Value = -61
$ws.Cells.Item(76, 8).Value = 2999.5
$ws.Cells.Item(76, 9).Value = 3000
$ws.Cells.Item(76, 10).Value = 2999.3333
$ws.Cells.Item(76, 11).Value = 3000
$ws.Cells.Item(76, 12).Value = 2999.3333
$ws.Cells.Item(76, 13).Value = -2685
$ws.Cells.Item(76, 14).Value = -3629.3333
$ws.Cells.Item(79, 8).Value = 2999.5
$ws.Cells.Item(79, 9).Value = 3000
$ws.Cells.Item(79, 10).Value = 2999.3333
$ws.Cells.Item(79, 11).Value = 3000
$ws.Cells.Item(79, 12).Value = 2999.3333
$ws.Cells.Item(79, 13).Value = -1908
$ws.Cells.Item(79, 14).Value = -5183.3333
$ws.Cells.Item(80, 8).Value = 1137.6666
$ws.Cells.Item(80, 9).Value = 672.25
$ws.Cells.Item(80, 10).Value = 2068.5
$ws.Cells.Item(80, 11).Value = 2016.75
$ws.Cells.Item(80, 12).Value = 6205.5
$ws.Cells.Item(80, 13).Value = -1018.75
$ws.Cells.Item(80, 14).Value = -8201.5
$ws.Cells.Item(83, 8).Value = 1137.6666
$ws.Cells.Item(83, 9).Value = 672.25
$ws.Cells.Item(83, 10).Value = 2068.5
$ws.Cells.Item(83, 11).Value = 6050.25
$ws.Cells.Item(83, 12).Value = 18616.5
$ws.Cells.Item(83, 13).Value = -1058.25
$ws.Cells.Item(83, 14).Value = -28600.5
$ws.Cells.Item(137, 8).Value = 2932.7222
$ws.Cells.Item(137, 9).Value = 2157.4167
$ws.Cells.Item(137, 11).Value = 6472.250100000001
$ws.Cells.Item(137, 13).Value = -3922.250100000001
$ws.Cells.Item(138, 8).Value = 3235.25
$ws.Cells.Item(138, 9).Value = 1470.5
$ws.Cells.Item(138, 10).Value = 5000
$ws.Cells.Item(138, 11).Value = 4411.5
$ws.Cells.Item(138, 12).Value = 15000
$ws.Cells.Item(138, 13).Value = 728.5
$ws.Cells.Item(138, 14).Value = -25280

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2310
$ws.Cells.Item(61, 9).Value = 2341
$ws.Cells.Item(61, 10).Value = 2155
$ws.Cells.Item(61, 11).Value = 2341
$ws.Cells.Item(61, 12).Value = 2155
$ws.Cells.Item(61, 13).Value = -2129
$ws.Cells.Item(61, 14).Value = -2579
$ws.Cells.Item(112, 8).Value = 39999
$ws.Cells.Item(112, 10).Value = 39999
$ws.Cells.Item(112, 12).Value = 39999
$ws.Cells.Item(112, 14).Value = -42953
$ws.Cells.Item(114, 8).Value = 0
$ws.Cells.Item(114, 10).Value = 0
$ws.Cells.Item(114, 14).Value = $null
$ws.Cells.Item(119, 8).Value = 55999.75
$ws.Cells.Item(119, 10).Value = 55999.75
$ws.Cells.Item(119, 12).Value = 55999.75
$ws.Cells.Item(119, 14).Value = -65675.75
$ws.Cells.Item(136, 8).Value = 2310
$ws.Cells.Item(136, 9).Value = 2341
$ws.Cells.Item(136, 10).Value = 2155
$ws.Cells.Item(136, 11).Value = 7023
$ws.Cells.Item(136, 12).Value = 6465
$ws.Cells.Item(136, 13).Value = -4473
$ws.Cells.Item(136, 14).Value = -11565

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(35, 8).Value = 14997.5
$ws.Cells.Item(35, 10).Value = 14997.5
$ws.Cells.Item(35, 12).Value = 14997.5
$ws.Cells.Item(35, 14).Value = -15617.5
$ws.Cells.Item(82, 8).Value = 23749.316
$ws.Cells.Item(82, 10).Value = 29998
$ws.Cells.Item(82, 12).Value = 29998
$ws.Cells.Item(82, 14).Value = -30764
$ws.Cells.Item(85, 8).Value = 23749.316
$ws.Cells.Item(85, 10).Value = 29998
$ws.Cells.Item(85, 12).Value = 29998
$ws.Cells.Item(85, 14).Value = -32650
$ws.Cells.Item(92, 8).Value = 37500
$ws.Cells.Item(92, 10).Value = 37500
$ws.Cells.Item(92, 12).Value = 37500
$ws.Cells.Item(92, 14).Value = -42492
$ws.Cells.Item(94, 8).Value = 1952
$ws.Cells.Item(94, 9).Value = 1776.2727
$ws.Cells.Item(94, 11).Value = 1776.2727
$ws.Cells.Item(94, 13).Value = -1325.2727
$ws.Cells.Item(107, 8).Value = 1078.2222
$ws.Cells.Item(107, 9).Value = 843.5714
$ws.Cells.Item(107, 11).Value = 843.5714
$ws.Cells.Item(107, 13).Value = 1076.4286

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(41, 8).Value = 11801.223
$ws.Cells.Item(41, 9).Value = 2070.3333
$ws.Cells.Item(41, 11).Value = 2070.3333
$ws.Cells.Item(41, 13).Value = -1642.3333
$ws.Cells.Item(50, 8).Value = 20102.3
$ws.Cells.Item(59, 8).Value = 28900.8
$ws.Cells.Item(60, 8).Value = 20155.857
$ws.Cells.Item(60, 10).Value = 20000
$ws.Cells.Item(60, 12).Value = 20000
$ws.Cells.Item(60, 14).Value = -21022
$ws.Cells.Item(74, 8).Value = 29714.285
$ws.Cells.Item(74, 10).Value = 29714.285
$ws.Cells.Item(74, 12).Value = 29714.285
$ws.Cells.Item(74, 14).Value = -31462.285
$ws.Cells.Item(77, 8).Value = 29714.285
$ws.Cells.Item(77, 10).Value = 29714.285
$ws.Cells.Item(77, 12).Value = 89142.855
$ws.Cells.Item(77, 14).Value = -97878.855
$ws.Cells.Item(105, 8).Value = 1049.75
$ws.Cells.Item(105, 9).Value = 1100
$ws.Cells.Item(105, 10).Value = 999.5
$ws.Cells.Item(105, 11).Value = 1100
$ws.Cells.Item(105, 12).Value = 999.5
$ws.Cells.Item(105, 13).Value = 647
$ws.Cells.Item(105, 14).Value = -4493.5
$ws.Cells.Item(107, 8).Value = 517.5294
$ws.Cells.Item(107, 9).Value = 475.16666
$ws.Cells.Item(107, 11).Value = 475.16666
$ws.Cells.Item(107, 13).Value = 1444.83334

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 1165.1
$ws.Cells.Item(68, 9).Value = 1062.3334
$ws.Cells.Item(68, 10).Value = 1209.1428
$ws.Cells.Item(68, 11).Value = 3187.0002
$ws.Cells.Item(68, 12).Value = 3627.4284
$ws.Cells.Item(68, 13).Value = -2376.0002
$ws.Cells.Item(68, 14).Value = -5249.428400000001
$ws.Cells.Item(71, 8).Value = 1165.1
$ws.Cells.Item(71, 9).Value = 1062.3334
$ws.Cells.Item(71, 10).Value = 1209.1428
$ws.Cells.Item(71, 11).Value = 9561.000599999999
$ws.Cells.Item(71, 12).Value = 10882.2852
$ws.Cells.Item(71, 13).Value = -5505.000599999999
$ws.Cells.Item(71, 14).Value = -18994.2852
$ws.Cells.Item(86, 8).Value = 1198.5
$ws.Cells.Item(86, 10).Value = 1198.5
$ws.Cells.Item(86, 12).Value = 3595.5
$ws.Cells.Item(86, 14).Value = -5967.5
$ws.Cells.Item(89, 8).Value = 1198.5
$ws.Cells.Item(89, 10).Value = 1198.5
$ws.Cells.Item(89, 12).Value = 10786.5
$ws.Cells.Item(89, 14).Value = -22642.5
$ws.Cells.Item(97, 8).Value = 1229.6666
$ws.Cells.Item(97, 10).Value = 1229.6666
$ws.Cells.Item(97, 12).Value = 3688.9998
$ws.Cells.Item(97, 14).Value = -4680.9998

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 628.1053000000001
$ws.Cells.Item(2, 9).Value = 489.58334
$ws.Cells.Item(2, 10).Value = 865.5714
$ws.Cells.Item(2, 11).Value = 489.58334
$ws.Cells.Item(2, 12).Value = 865.5714
$ws.Cells.Item(2, 13).Value = -376.58334
$ws.Cells.Item(2, 14).Value = -1091.5714

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 2135.5454
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 14).Value = $null

Write-Output "Applied all changes"